# Commit: table style swap on the Sources-of-finance table (slide 6) plus a
# theme palette swap (deck's theme goes from the "Integral" palette to the
# stock "Office Theme" palette).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
# Slide 6 has a single table shape ("Google Shape;127;p18"); find it by
# scanning for HasTable rather than hard-coding the shape index.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{B37493BE-E322-468D-BB73-72A862BC0E1A}")
        }
    }
}

# --- 2. Theme colors: Integral palette -> Office Theme palette ------------
# ThemeColorScheme.Colors(n) order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10
# accent1..accent6, 11 hlink, 12 folHlink. RGB values are the usual packed
# R + G*256 + B*65536 integer (e.g. 0xBBGGRR as a literal).
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
